# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for Chirimoya (Cultivar IV Region) right
# before the existing row 411, shifting all subsequent rows down by 3
# (dimension grows from A1:T504 to A1:T507).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 411..413 (pushes old 411.. down to 414..)
$ws.Range("A411:A413").EntireRow.Insert()

# Shared values for the three new rows (same market/product as the rest
# of the sheet).
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$fecha       = 45218
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100107
$producto    = "Otros"
$categoriaId = 100107002
$categoria   = "Chirimoya"
$variedad    = "Cultivar IV Región"
$unidad      = "$/bandeja 10 kilos"
$origen      = "Provincia de Limarí"

function Set-ChirimoyaRow($row, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg, $kgUnidad) {
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-ChirimoyaRow 411 "Especial" 200 26000 26000 26000 2600 10
Set-ChirimoyaRow 412 "Primera"  300 23000 23000 23000 2300 10
Set-ChirimoyaRow 413 "Segunda"  250 20000 20000 20000 2000 10

Write-Output "done"
